# aggiornamenti ad A* e a tutti i file per l'esame di Settembre
#
# Applies the data + view-state edits described by the commit:
#  - Sheet "Mondo blocchi 1 - Profondita su" (sheet1): A* column (G2) time value updated
#  - Sheet "Labirinto 10x10 - Profondita su" (sheet3): A* column values updated (time/steps rows)
#  - Sheet "Labirinto 20x20 - Profondita su" (sheet4): A* column (G2) time value updated
#  - Active sheet moves from "Labirinto 10x10" to "Mondo blocchi 2"
#  - Selected cell on a couple of sheets changes to track the edited column

$wb = $excel.ActiveWorkbook

$wsMondo1     = $wb.Worksheets.Item(1)   # Mondo blocchi 1 - Profondita su
$wsMondo2     = $wb.Worksheets.Item(2)   # Mondo blocchi 2 - Profondita su
$wsLab10      = $wb.Worksheets.Item(3)   # Labirinto 10x10 - Profondita su
$wsLab20      = $wb.Worksheets.Item(4)   # Labirinto 20x20 - Profondita su

# --- Mondo blocchi 1: A* time (column G, row 2) ---
$wsMondo1.Activate()
$wsMondo1.Range("G2").Value = 70
$wsMondo1.Range("G3").Select()

# --- Labirinto 10x10: A* column (row2 = tempo, row3 = numero di passi) ---
$wsLab10.Activate()
$wsLab10.Range("D2").Value = 190
$wsLab10.Range("F2").Value = 1217
$wsLab10.Range("G2").Value = 6
$wsLab10.Range("H2").Value = 272

$wsLab10.Range("B3").Value = 49
$wsLab10.Range("D3").Value = 19
$wsLab10.Range("E3").Value = 19
$wsLab10.Range("F3").Value = 19
$wsLab10.Range("G3").Value = 19
$wsLab10.Range("H3").Value = 19
$wsLab10.Range("H3").Select()

# --- Labirinto 20x20: A* time (column G, row 2) ---
$wsLab20.Activate()
$wsLab20.Range("G2").Value = 15
$wsLab20.Range("G2").Select()

# --- Finally, make "Mondo blocchi 2" the active sheet/tab ---
$wsMondo2.Activate()
$wsMondo2.Range("E2").Select()
